$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '51.291.35'
$ws.Range("E2").Value = '  -0.07%  '
$ws.Range("D3").Value = '2.974.51'
$ws.Range("E3").Value = '  +1.75%  '
$ws.Range("E4").Value = '  +0.29%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '380.84'
$ws.Range("E5").Value = '  +1.92%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '103.13'
$ws.Range("E6").Value = '  -1.16%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.542'
$ws.Range("E7").Value = '  -0.84%  '
$ws.Range("E8").Value = '  +0.04%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.590'
$ws.Range("E9").Value = '  -0.32%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '36.85'
$ws.Range("E10").Value = '  -1.17%  '
$ws.Range("E11").Value = '  -0.42%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.0844'
$ws.Range("E12").Value = '  +0.34%  '
$ws.Range("D13").Value = '3.443.45'
$ws.Range("E13").Value = '  +1.89%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '18.19'
$ws.Range("E14").Value = '  -1.17%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '7.50'
$ws.Range("E15").Value = '  +1.35%  '
$ws.Range("D16").Value = '2.976.06'
$ws.Range("E16").Value = '  +2.02%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.996'
$ws.Range("E17").Value = '  +6.26%  '
$ws.Range("D18").Value = '51.312.04'
$ws.Range("E18").Value = '  +0.09%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '3.25'
$ws.Range("E19").Value = '  -2.48%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '7.37'
$ws.Range("E20").Value = '  +0.82%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '12.75'
$ws.Range("E21").Value = '  -2.41%  '
$ws.Range("D22").Value = '0.0₃0959'
$ws.Range("E22").Value = '  +1.21%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '68.99'
$ws.Range("E23").Value = '  +0.64%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '262.36'
$ws.Range("E24").Value = '  +0.37%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.90'
$ws.Range("E25").Value = '  +7.29%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '8.15'
$ws.Range("E26").Value = '  +13.17%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '7.57'
$ws.Range("E27").Value = '  +11.31%  '
$ws.Range("E28").Value = '  +15.06%  '
$ws.Range("E29").Value = '  -1.32%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '4.12'
$ws.Range("E30").Value = '  -0.17%  '
$ws.Range("E31").Value = '  -0.03%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '25.90'
$ws.Range("E32").Value = '  +0.05%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '9.82'
$ws.Range("E33").Value = '  -1.24%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '34.29'
$ws.Range("E34").Value = '  -1.29%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '50.91'
$ws.Range("E35").Value = '  -0.25%  '
$ws.Range("E36").Value = '  -2.31%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.0449'
$ws.Range("E37").Value = '  +5.56%  '
$ws.Range("E38").Value = '  +0.23%  '
$ws.Range("E39").Value = '  -1.56%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '16.86'
$ws.Range("E40").Value = '  -1.77%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '2.56'
$ws.Range("E41").Value = '  -0.69%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.115'
$ws.Range("E42").Value = '  +1.57%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.82'
$ws.Range("E43").Value = '  -1.77%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '122.43'
$ws.Range("E44").Value = '  +2.39%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '21.59'
$ws.Range("E45").Value = '  -2.49%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.06'
$ws.Range("E46").Value = '  -1.81%  '
$ws.Range("E47").Value = '  +7.65%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '3.29'
$ws.Range("E49").Value = '  +2.73%  '
$ws.Range("D50").Value = '2.028.75'
$ws.Range("E50").Value = '  -0.05%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0331'
$ws.Range("E51").Value = '  +3.12%  '
